$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.738.09"
$ws.Range("E2").Value = "  -0.34%  "

$ws.Range("D3").Value = "2.776.30"
$ws.Range("E3").Value = "  -1.51%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "356.76"
$ws.Range("E5").Value = "  +0.27%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.29"
$ws.Range("E6").Value = "  -1.94%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.556"
$ws.Range("E7").Value = "  -1.32%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.590"
$ws.Range("E9").Value = "  -1.36%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.80"
$ws.Range("E10").Value = "  -2.43%  "

$ws.Range("E11").Value = "  +2.44%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0843"
$ws.Range("E12").Value = "  -1.39%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.46"
$ws.Range("E13").Value = "  -2.27%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.60"
$ws.Range("E14").Value = "  -2.75%  "

$ws.Range("D15").Value = "3.217.95"
$ws.Range("E15").Value = "  -1.40%  "

$ws.Range("D16").Value = "2.782.23"
$ws.Range("E16").Value = "  -1.12%  "

$ws.Range("E17").Value = "  +1.80%  "

$ws.Range("D18").Value = "51.708.14"
$ws.Range("E18").Value = "  -0.28%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.43"
$ws.Range("E19").Value = "  -1.56%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.00"
$ws.Range("E20").Value = "  -3.87%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.14"
$ws.Range("E21").Value = "  -1.87%  "

$ws.Range("D22").Value = "0.0₃0969"
$ws.Range("E22").Value = "  -2.44%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.14"
$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "269.05"
$ws.Range("E24").Value = "  +0.47%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.73"
$ws.Range("E25").Value = "  -2.73%  "

$ws.Range("E26").Value = "  -2.28%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.09%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.163"
$ws.Range("E28").Value = "  +16.17%  "

$ws.Range("E29").Value = "  -0.61%  "

$ws.Range("E30").Value = "  -1.52%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.27"
$ws.Range("E31").Value = "  +6.42%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "51.88"
$ws.Range("E32").Value = "  -1.25%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.61"
$ws.Range("E33").Value = "  -0.30%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0448"
$ws.Range("E34").Value = "  -12.72%  "

$ws.Range("E35").Value = "  -0.48%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.12"
$ws.Range("E36").Value = "  -6.33%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.74"
$ws.Range("E38").Value = "  +2.31%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.14"
$ws.Range("E39").Value = "  -4.59%  "

$ws.Range("E40").Value = "  -3.77%  "

$ws.Range("E41").Value = "  +0.61%  "

$ws.Range("E42").Value = "  -2.24%  "

$ws.Range("E43").Value = "  -1.76%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "119.67"
$ws.Range("E44").Value = "  -5.18%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.64"
$ws.Range("E45").Value = "  -7.38%  "

$ws.Range("D46").Value = "2.079.78"
$ws.Range("E46").Value = "  -0.88%  "

$ws.Range("E47").Value = "  -2.33%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.28"
$ws.Range("E48").Value = "  +0.80%  "

$ws.Range("E49").Value = "  -5.00%  "

$ws.Range("E50").Value = "  -5.11%  "

$ws.Range("E51").Value = "  +0.86%  "
